# Auto-generated: bulk market-price data refresh across all item sheets
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2697.5
$ws.Range("I40").Value = 2466.6667
$ws.Range("J40").Value = 2836
$ws.Range("K40").Value = 2466.6667
$ws.Range("L40").Value = 2836
$ws.Range("M40").Value = -2291.6667
$ws.Range("N40").Value = -3186
$ws.Range("H64").Value = 3973.658
$ws.Range("I64").Value = 3836.3635
$ws.Range("J64").Value = 4029.5925
$ws.Range("K64").Value = 3836.3635
$ws.Range("L64").Value = 4029.5925
$ws.Range("M64").Value = -3588.3635
$ws.Range("N64").Value = -4525.592500000001
$ws.Range("H67").Value = 3973.658
$ws.Range("I67").Value = 3836.3635
$ws.Range("J67").Value = 4029.5925
$ws.Range("K67").Value = 3836.3635
$ws.Range("L67").Value = 4029.5925
$ws.Range("M67").Value = -2978.3635
$ws.Range("N67").Value = -5745.592500000001
$ws.Range("H74").Value = 3260.2
$ws.Range("I74").Value = 3003
$ws.Range("J74").Value = 3278.5715
$ws.Range("K74").Value = 3003
$ws.Range("L74").Value = 3278.5715
$ws.Range("M74").Value = -2067
$ws.Range("N74").Value = -5150.5715
$ws.Range("H76").Value = 3550.0264
$ws.Range("J76").Value = 3550.0264
$ws.Range("L76").Value = 3550.0264
$ws.Range("N76").Value = -4180.026400000001
$ws.Range("H77").Value = 3260.2
$ws.Range("I77").Value = 3003
$ws.Range("J77").Value = 3278.5715
$ws.Range("K77").Value = 15015
$ws.Range("L77").Value = 16392.8575
$ws.Range("M77").Value = -10335
$ws.Range("N77").Value = -25752.8575
$ws.Range("H79").Value = 3550.0264
$ws.Range("J79").Value = 3550.0264
$ws.Range("L79").Value = 3550.0264
$ws.Range("N79").Value = -5734.026400000001
$ws.Range("H113").Value = 1762.1875
$ws.Range("J113").Value = 1850.0454
$ws.Range("L113").Value = 1850.0454
$ws.Range("N113").Value = -8358.045399999999
$ws.Range("H116").Value = 2334036.2
$ws.Range("I116").Value = 6995945
$ws.Range("J116").Value = 3081.818
$ws.Range("K116").Value = 6995945
$ws.Range("L116").Value = 3081.818
$ws.Range("M116").Value = -6992503
$ws.Range("N116").Value = -9965.817999999999
$ws.Range("H121").Value = 3994.5
$ws.Range("J121").Value = 3994.5
$ws.Range("L121").Value = 11983.5
$ws.Range("N121").Value = -15477.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 26259.25
$ws.Range("I74").Value = 26259.25
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 26259.25
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -25385.25
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 26259.25
$ws.Range("I77").Value = 26259.25
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 131296.25
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -126928.25
$ws.Range("N77").Value = ""
$ws.Range("H122").Value = 1476.3043
$ws.Range("I122").Value = 1428.25
$ws.Range("J122").Value = 1796.6666
$ws.Range("K122").Value = 4284.75
$ws.Range("L122").Value = 5389.9998
$ws.Range("M122").Value = -1834.75
$ws.Range("N122").Value = -10289.9998
$ws.Range("H131").Value = 45655.715
$ws.Range("J131").Value = 45655.715
$ws.Range("L131").Value = 45655.715
$ws.Range("N131").Value = -55735.715
$ws.Range("H134").Value = 32429
$ws.Range("J134").Value = 32429
$ws.Range("L134").Value = 32429
$ws.Range("N134").Value = -42569
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2138.5151
$ws.Range("I20").Value = 1310.2778
$ws.Range("J20").Value = 3132.4
$ws.Range("K20").Value = 1310.2778
$ws.Range("L20").Value = 3132.4
$ws.Range("M20").Value = -1063.2778
$ws.Range("N20").Value = -3626.4
$ws.Range("H80").Value = 122.3
$ws.Range("I80").Value = 188.75
$ws.Range("J80").Value = 78
$ws.Range("K80").Value = 188.75
$ws.Range("L80").Value = 78
$ws.Range("M80").Value = 809.25
$ws.Range("N80").Value = -2074
$ws.Range("H83").Value = 122.3
$ws.Range("I83").Value = 188.75
$ws.Range("J83").Value = 78
$ws.Range("K83").Value = 943.75
$ws.Range("L83").Value = 390
$ws.Range("M83").Value = 4048.25
$ws.Range("N83").Value = -10374
$ws.Range("H94").Value = 606.7241
$ws.Range("I94").Value = 626.3043
$ws.Range("J94").Value = 531.6667
$ws.Range("K94").Value = 626.3043
$ws.Range("L94").Value = 531.6667
$ws.Range("M94").Value = -175.3043
$ws.Range("N94").Value = -1433.6667
$ws.Range("H134").Value = 2455.3684
$ws.Range("I134").Value = 2088.4285
$ws.Range("J134").Value = 3482.8
$ws.Range("K134").Value = 6265.2855
$ws.Range("L134").Value = 10448.4
$ws.Range("M134").Value = -3730.2855
$ws.Range("N134").Value = -15518.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 527
$ws.Range("I107").Value = 418.14285
$ws.Range("J107").Value = 831.8
$ws.Range("K107").Value = 418.14285
$ws.Range("L107").Value = 831.8
$ws.Range("M107").Value = 1501.85715
$ws.Range("N107").Value = -4671.8
$ws.Range("H122").Value = 758.14813
$ws.Range("J122").Value = 553.2
$ws.Range("L122").Value = 1659.6
$ws.Range("N122").Value = -6559.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 459.6
$ws.Range("I23").Value = 399.33334
$ws.Range("J23").Value = 550
$ws.Range("K23").Value = 1198.00002
$ws.Range("L23").Value = 1650
$ws.Range("M23").Value = -963.0000199999999
$ws.Range("N23").Value = -2120
$ws.Range("H97").Value = 532.5625
$ws.Range("I97").Value = 329
$ws.Range("J97").Value = 654.7
$ws.Range("K97").Value = 987
$ws.Range("L97").Value = 1964.1
$ws.Range("M97").Value = -491
$ws.Range("N97").Value = -2956.1
$ws.Range("H107").Value = 503.85715
$ws.Range("I107").Value = 127.7
$ws.Range("J107").Value = 845.8182
$ws.Range("K107").Value = 383.1
$ws.Range("L107").Value = 2537.4546
$ws.Range("M107").Value = 1536.9
$ws.Range("N107").Value = -6377.4546
$ws.Range("H122").Value = 6508.6665
$ws.Range("I122").Value = 8930.120000000001
$ws.Range("J122").Value = 1005.36365
$ws.Range("K122").Value = 80371.08
$ws.Range("L122").Value = 9048.272849999999
$ws.Range("M122").Value = -77921.08
$ws.Range("N122").Value = -13948.27285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8598.462
$ws.Range("I70").Value = 8598.462
$ws.Range("K70").Value = 8598.462
$ws.Range("M70").Value = -8328.462
$ws.Range("H73").Value = 8598.462
$ws.Range("I73").Value = 8598.462
$ws.Range("K73").Value = 8598.462
$ws.Range("M73").Value = -7662.462
$ws.Range("H93").Value = 50327
$ws.Range("J93").Value = 50327
$ws.Range("L93").Value = 50327
$ws.Range("N93").Value = -54071
$ws.Range("H122").Value = 14287094
$ws.Range("I122").Value = 20001290
$ws.Range("K122").Value = 60003870
$ws.Range("M122").Value = -60001420
$ws.Range("H126").Value = 2973.7568
$ws.Range("I126").Value = 2667.0833
$ws.Range("J126").Value = 3539.923
$ws.Range("K126").Value = 8001.249899999999
$ws.Range("L126").Value = 10619.769
$ws.Range("M126").Value = -5531.249899999999
$ws.Range("N126").Value = -15559.769
$ws.Range("H132").Value = 4484.55
$ws.Range("I132").Value = 4277.2354
$ws.Range("K132").Value = 12831.7062
$ws.Range("M132").Value = -10301.7062
$ws.Range("H141").Value = 33000
$ws.Range("J141").Value = 33000
$ws.Range("L141").Value = 33000
$ws.Range("N141").Value = -43360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4725
$ws.Range("I122").Value = 2320
$ws.Range("J122").Value = 8733.333000000001
$ws.Range("K122").Value = 6960
$ws.Range("L122").Value = 26199.999
$ws.Range("M122").Value = -4510
$ws.Range("N122").Value = -31099.999
$ws.Range("H135").Value = 28429
$ws.Range("J135").Value = 28429
$ws.Range("L135").Value = 28429
$ws.Range("N135").Value = -38569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43286
$ws.Range("J46").Value = 43286
$ws.Range("L46").Value = 43286
$ws.Range("N46").Value = -43748
$ws.Range("H107").Value = 492.1111
$ws.Range("I107").Value = 461.2857
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1383.8571
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 536.1428999999998
$ws.Range("N107").Value = -5640
$ws.Range("H122").Value = 2153.7354
$ws.Range("I122").Value = 1995.9
$ws.Range("J122").Value = 2379.2144
$ws.Range("K122").Value = 5987.700000000001
$ws.Range("L122").Value = 7137.6432
$ws.Range("M122").Value = -3537.700000000001
$ws.Range("N122").Value = -12037.6432
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").Value = ""
$ws.Range("H134").Value = 43286
$ws.Range("J134").Value = 43286
$ws.Range("L134").Value = 129858
$ws.Range("N134").Value = -134928

Write-Host "Updated market price data across all sheets"
